$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 5867
$ws1.Range("F5").Value = 5867
$ws1.Range("F7").Value = 2910
$ws1.Range("F8").Value = 1254
$ws1.Range("F9").Value = 389
$ws1.Range("F13").Value = 673
$ws1.Range("F14").Value = 164
$ws1.Range("F15").Value = 4182
$ws1.Range("F16").Value = 4182
$ws1.Range("F18").Value = 82
$ws1.Range("F21").Value = 189
$ws1.Range("F23").Value = 6270
$ws1.Range("F24").Value = 6270
$ws1.Range("F25").Value = 222
$ws1.Range("F28").Value = 427
$ws1.Range("F29").Value = 204
$ws1.Range("F32").Value = 1604
$ws1.Range("F34").Value = 1739
$ws1.Range("F35").Value = 5839
$ws1.Range("F36").Value = 95
$ws1.Range("F39").Value = 70
$ws1.Range("F40").Value = 121
$ws1.Range("F41").Value = 3911
$ws1.Range("F42").Value = 121
$ws1.Range("F43").Value = 73
$ws1.Range("F45").Value = 2386
$ws1.Range("F48").Value = 1000
$ws1.Range("F50").Value = 262

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 183
$ws2.Range("F11").Value = 14

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 5867
$ws4.Range("F5").Value = 5867
$ws4.Range("F7").Value = 2910
$ws4.Range("F8").Value = 1254
$ws4.Range("F12").Value = 183
$ws4.Range("F13").Value = 673
$ws4.Range("F14").Value = 164
$ws4.Range("F15").Value = 4182
$ws4.Range("F16").Value = 4182
$ws4.Range("F18").Value = 82
$ws4.Range("F21").Value = 189
$ws4.Range("F23").Value = 6270
$ws4.Range("F24").Value = 6270
$ws4.Range("F25").Value = 222
$ws4.Range("F27").Value = 427
$ws4.Range("F28").Value = 204
$ws4.Range("F32").Value = 1604
$ws4.Range("F35").Value = 1739
$ws4.Range("F37").Value = 5839
$ws4.Range("F38").Value = 95
$ws4.Range("F40").Value = 70
$ws4.Range("F41").Value = 3912
$ws4.Range("F42").Value = 73
$ws4.Range("F46").Value = 2386
$ws4.Range("F49").Value = 1000
$ws4.Range("F50").Value = 262
$ws4.Range("F52").Value = 14

